$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in D19: "8.30 to4.45" -> "8.30 to 4.45"
$ws.Range("D19").Value = "8.30 to 4.45"

# Add new row 20 with the next day's entry
$ws.Range("A20").Value = 15
$ws.Range("B20").Value = "Prabha"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "04.01.2018"
$ws.Range("E20").Value = "searched videos for API connection,json"
$ws.Range("D20").Value = "8.40 to 4.45"
$ws.Range("F20").Value = "completed"

# Copy formatting from row 19 to row 20 so the new row matches the table style
$ws.Range("A19:F19").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(20).RowHeight = $ws.Rows.Item(19).RowHeight

# Update the selection to match the new active cell
$ws.Range("E20").Select()
